$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename fund / portfolio companies ---
$ws.Range("A2:A5").Value = "Demo Fund 2"
$ws.Range("B2").Value = "TSTF2 Port Co 3"
$ws.Range("B3").Value = "TSTF2 Port Co 3"
$ws.Range("B4").Value = "TSTF2 Port Co 4"
$ws.Range("B5").Value = "TSTF2 Port Co 3"

# --- Update investment dates ---
$ws.Range("C4").Value = 45214
$ws.Range("C5").Value = 45717

# --- Row 5 instrument changes from CCPS to Equity ---
$ws.Range("G5").Value = "Equity"

# --- Update investment amounts (recalculates E formulas in place) ---
$ws.Range("D2").Value = 100000000
$ws.Range("D3").Value = 40000000
$ws.Range("D4").Value = 280000000

# --- Row 5 ratio becomes a plain (non-formula) value ---
$ws.Range("E5").Value = -500000
$ws.Range("D5").Value = 100000000

# --- Drop the stray empty cells trailing row 5 ---
$ws.Range("J5").Clear()
$ws.Range("K5").Clear()

# --- Re-knit the shared formula group so it spans E2:E4 ---
$ws.Range("E2:E4").Formula = "=D2/F2"
$ws.Range("B2").Copy()
$ws.Range("E2:E4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Drop the now-unused trailing blank row ---
$ws.Rows(6).Delete()

# --- Shrink the hidden filter-database name to match the trimmed data ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$ALZ`$4"
    }
}

# --- Match the saved selection ---
$ws.Range("G5").Select()
